$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1285.4
$ws.Range("I107").Value = 1313.4615
$ws.Range("J107").Value = 1103
$ws.Range("K107").Value = 1313.4615
$ws.Range("L107").Value = 1103
$ws.Range("M107").Value = 606.5385000000001
$ws.Range("N107").Value = -4943

$ws.Range("H113").Value = 4006.7856
$ws.Range("I113").Value = 3964.1428
$ws.Range("J113").Value = 4049.4285
$ws.Range("K113").Value = 3964.1428
$ws.Range("L113").Value = 4049.4285
$ws.Range("M113").Value = -710.1428000000001
$ws.Range("N113").Value = -10557.4285

$ws.Range("H116").Value = 3903.0435
$ws.Range("I116").Value = 3220.9092
$ws.Range("J116").Value = 4528.3335
$ws.Range("K116").Value = 3220.9092
$ws.Range("L116").Value = 4528.3335
$ws.Range("M116").Value = 221.0907999999999
$ws.Range("N116").Value = -11412.3335

$ws.Range("H132").Value = 28582286
$ws.Range("I132").Value = 28582286
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 85746858
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -85744328
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 6857.9
$ws.Range("I137").Value = 12277
$ws.Range("J137").Value = 3245.1667
$ws.Range("K137").Value = 36831
$ws.Range("L137").Value = 9735.500100000001
$ws.Range("M137").Value = -34281
$ws.Range("N137").Value = -14835.5001

$ws.Range("H138").Value = 5061.8506
$ws.Range("I138").Value = 2612.9092
$ws.Range("J138").Value = 5890.723
$ws.Range("K138").Value = 7838.7276
$ws.Range("L138").Value = 17672.169
$ws.Range("M138").Value = -2698.7276
$ws.Range("N138").Value = -27952.169

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 50003064
$ws.Range("I2").Value = 75002390
$ws.Range("J2").Value = 4400
$ws.Range("K2").Value = 75002390
$ws.Range("L2").Value = 4400
$ws.Range("M2").Value = -75002277
$ws.Range("N2").Value = -4626

$ws.Range("H63").Value = 2328.4285
$ws.Range("J63").Value = 2274.75
$ws.Range("L63").Value = 2274.75
$ws.Range("N63").Value = -3646.75

$ws.Range("H66").Value = 2328.4285
$ws.Range("J66").Value = 2274.75
$ws.Range("L66").Value = 11373.75
$ws.Range("N66").Value = -18237.75

$ws.Range("H74").Value = 2434.56
$ws.Range("I74").Value = 1819.125
$ws.Range("K74").Value = 1819.125
$ws.Range("M74").Value = -945.125

$ws.Range("H77").Value = 2434.56
$ws.Range("I77").Value = 1819.125
$ws.Range("K77").Value = 9095.625
$ws.Range("M77").Value = -4727.625

$ws.Range("H113").Value = 30132.666
$ws.Range("J113").Value = 30132.666
$ws.Range("L113").Value = 30132.666
$ws.Range("N113").Value = -38810.666

$ws.Range("H116").Value = 50003064
$ws.Range("I116").Value = 75002390
$ws.Range("J116").Value = 4400
$ws.Range("K116").Value = 75002390
$ws.Range("L116").Value = 4400
$ws.Range("M116").Value = -75000096
$ws.Range("N116").Value = -8988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 50003064
$ws.Range("I3").Value = 75002390
$ws.Range("J3").Value = 4400
$ws.Range("K3").Value = 75002390
$ws.Range("L3").Value = 4400
$ws.Range("M3").Value = -75002276
$ws.Range("N3").Value = -4628

$ws.Range("H99").Value = 3270.5
$ws.Range("I99").Value = 2819
$ws.Range("J99").Value = 4625
$ws.Range("K99").Value = 2819
$ws.Range("L99").Value = 4625
$ws.Range("M99").Value = -1321
$ws.Range("N99").Value = -7621

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1292.3334
$ws.Range("I16").Value = 747.0769
$ws.Range("K16").Value = 747.0769
$ws.Range("M16").Value = -460.0769

$ws.Range("H18").Value = 32000
$ws.Range("J18").Value = 32000
$ws.Range("L18").Value = 32000
$ws.Range("N18").Value = -32460

$ws.Range("H31").Value = 3205.2163
$ws.Range("I31").Value = 2169.175
$ws.Range("J31").Value = 4424.0884
$ws.Range("K31").Value = 2169.175
$ws.Range("L31").Value = 4424.0884
$ws.Range("M31").Value = -1874.175
$ws.Range("N31").Value = -5014.0884

$ws.Range("H34").Value = 3205.2163
$ws.Range("I34").Value = 2169.175
$ws.Range("J34").Value = 4424.0884
$ws.Range("K34").Value = 2169.175
$ws.Range("L34").Value = 4424.0884
$ws.Range("M34").Value = -1967.175
$ws.Range("N34").Value = -4828.0884

$ws.Range("H58").Value = 7044981
$ws.Range("I58").Value = 1494.2325
$ws.Range("J58").Value = 17861764
$ws.Range("K58").Value = 1494.2325
$ws.Range("L58").Value = 17861764
$ws.Range("M58").Value = -1291.2325
$ws.Range("N58").Value = -17862170

$ws.Range("H113").Value = 1292.3334
$ws.Range("I113").Value = 747.0769
$ws.Range("K113").Value = 747.0769
$ws.Range("M113").Value = 1422.9231

$ws.Range("H136").Value = 7044981
$ws.Range("I136").Value = 1494.2325
$ws.Range("J136").Value = 17861764
$ws.Range("K136").Value = 4482.6975
$ws.Range("L136").Value = 53585292
$ws.Range("M136").Value = -1932.6975
$ws.Range("N136").Value = -53590392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10880.462
$ws.Range("I87").Value = 3942.8
$ws.Range("J87").Value = 15216.5
$ws.Range("K87").Value = 11828.4
$ws.Range("L87").Value = 45649.5
$ws.Range("M87").Value = -10580.4
$ws.Range("N87").Value = -48145.5

$ws.Range("H90").Value = 10880.462
$ws.Range("I90").Value = 3942.8
$ws.Range("J90").Value = 15216.5
$ws.Range("K90").Value = 35485.2
$ws.Range("L90").Value = 136948.5
$ws.Range("M90").Value = -29245.2
$ws.Range("N90").Value = -149428.5

$ws.Range("H120").Value = 16530.6
$ws.Range("I120").Value = 4391.8
$ws.Range("J120").Value = 22600
$ws.Range("K120").Value = 13175.4
$ws.Range("L120").Value = 67800
$ws.Range("M120").Value = -8337.400000000001
$ws.Range("N120").Value = -77476

$ws.Range("H132").Value = 2715.4736
$ws.Range("I132").Value = 1650.5
$ws.Range("J132").Value = 3490
$ws.Range("K132").Value = 14854.5
$ws.Range("L132").Value = 31410
$ws.Range("M132").Value = -12324.5
$ws.Range("N132").Value = -36470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23751.75
$ws.Range("I15").Value = 12007
$ws.Range("J15").Value = 27666.666
$ws.Range("K15").Value = 12007
$ws.Range("L15").Value = 27666.666
$ws.Range("M15").Value = -11719
$ws.Range("N15").Value = -28242.666

$ws.Range("H81").Value = 23751.75
$ws.Range("I81").Value = 12007
$ws.Range("J81").Value = 27666.666
$ws.Range("K81").Value = 12007
$ws.Range("L81").Value = 27666.666
$ws.Range("M81").Value = -11009
$ws.Range("N81").Value = -29662.666

$ws.Range("H84").Value = 23751.75
$ws.Range("I84").Value = 12007
$ws.Range("J84").Value = 27666.666
$ws.Range("K84").Value = 36021
$ws.Range("L84").Value = 82999.99800000001
$ws.Range("M84").Value = -31029
$ws.Range("N84").Value = -92983.99800000001

$ws.Range("H102").Value = 75777.07000000001
$ws.Range("I102").Value = 3432.2222
$ws.Range("J102").Value = 205997.8
$ws.Range("K102").Value = 3432.2222
$ws.Range("L102").Value = 205997.8
$ws.Range("M102").Value = -1810.2222
$ws.Range("N102").Value = -209241.8

$ws.Range("H126").Value = 1254314.2
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030

$ws.Range("H132").Value = 4658.324
$ws.Range("I132").Value = 5725.5557
$ws.Range("J132").Value = 3647.2632
$ws.Range("K132").Value = 17176.6671
$ws.Range("L132").Value = 10941.7896
$ws.Range("M132").Value = -14646.6671
$ws.Range("N132").Value = -16001.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 29166.666
$ws.Range("J106").Value = 29166.666
$ws.Range("L106").Value = 29166.666
$ws.Range("N106").Value = -31690.666

$ws.Range("H132").Value = 2201.9019
$ws.Range("I132").Value = 1619.3611
$ws.Range("K132").Value = 4858.0833
$ws.Range("M132").Value = -2328.0833

$ws.Range("H135").Value = 40937.5
$ws.Range("J135").Value = 40937.5
$ws.Range("L135").Value = 40937.5
$ws.Range("N135").Value = -51077.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 854.4
$ws.Range("I3").Value = 560
$ws.Range("J3").Value = 928
$ws.Range("K3").Value = 560
$ws.Range("L3").Value = 928
$ws.Range("M3").Value = -446
$ws.Range("N3").Value = -1156

$ws.Range("H8").Value = 351.5
$ws.Range("I8").Value = 203
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 203
$ws.Range("L8").Value = 500
$ws.Range("M8").Value = -63
$ws.Range("N8").Value = -780

$ws.Range("H92").Value = 29950
$ws.Range("J92").Value = 29950
$ws.Range("L92").Value = 29950
$ws.Range("N92").Value = -34942
